# Add a "description" header column to the "params" sheet (M1), copy the
# header style from the existing header cells, and move the active
# selection to M1 to match the post-edit workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# New header cell with the same style as the other header cells (e.g. L1).
$ws.Range("M1").Value = "description"
$ws.Range("M1").Style = $ws.Range("L1").Style

# Make this sheet the active one and move the selection to the new cell,
# matching the workbook's recorded UI state after the edit.
$ws.Activate()
$ws.Range("M1").Select()
